$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: new blood-pressure / heart-rate / SpO2 reading ---
$ws.Range("A19").Value = 45156.743055555555
$ws.Range("B19").Formula = "=(134+128)/2"
$ws.Range("C19").Formula = "=(94+89)/2"
$ws.Range("D19").Formula = "=(85+78)/2"
$ws.Range("E19").Value = 98

# --- Row 20: new blood-pressure / heart-rate / SpO2 reading ---
$ws.Range("A20").Value = 45157.40625
$ws.Range("B20").Formula = "=(129+136)/2"
$ws.Range("C20").Formula = "=(91+93)/2"
$ws.Range("D20").Formula = "=(97+103)/2"
$ws.Range("E20").Value = 96

# --- Update the active selection to reflect where the user left off ---
$ws.Range("C17").Select()
